# Updates the cryptos price/volume table to the latest scraped snapshot.
# Price cells (column D) are leading-apostrophe prefixed so Excel stores
# them as literal text (matching the original inlineStr cells) instead of
# auto-coercing numeric-looking strings like "5.280" into the number 5.28.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.365.37"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "'1.794.66"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'306.94"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("D7").Value = "'0.4514"
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("D8").Value = "'0.3597"
$ws.Range("E8").Value = "  -2.80%  "
$ws.Range("D9").Value = "'45.88"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "'0.07077"
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").Value = "'0.8831"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "'0.07753"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "'19.46"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "'1.793.35"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").Value = "'5.280"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "'6.319"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "'84.84"
$ws.Range("E17").Value = "  -2.81%  "
$ws.Range("D18").Value = "'1.007"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'0.000008508"
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D21").Value = "'14.25"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").Value = "'26.387.46"
$ws.Range("E22").Value = "  -2.20%  "
$ws.Range("D23").Value = "'4.968"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "'2.039.19"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "'1.974"
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("D27").Value = "'151.03"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").Value = "'17.81"
$ws.Range("E28").Value = "  -2.27%  "
$ws.Range("D29").Value = "'2.013"
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("D30").Value = "'111.85"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("D31").Value = "'4.869"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").Value = "'0.08666"
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").Value = "'3.069"
$ws.Range("E33").Value = "  +2.62%  "
$ws.Range("D34").Value = "'4.439"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7225"
$ws.Range("E35").Value = "  -3.72%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'2.710"
$ws.Range("E36").Value = "  +5.22%  "
$ws.Range("D37").Value = "'1.105"
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").Value = "'1.066"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("D40").Value = "'0.01930"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").Value = "'0.05086"
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("D42").Value = "'2.861"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'6.842"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.5040"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("D45").Value = "'0.1513"
$ws.Range("E45").Value = "  -5.48%  "
$ws.Range("D46").Value = "'7.997"
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").Value = "'0.4619"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("D49").Value = "'101.14"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("D50").Value = "'9.810"
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("D51").Value = "'1.575"
$ws.Range("E51").Value = "  -2.35%  "
